$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point precision difference on A8 (timestamp)
$ws.Cells.Item(8, 1).Value = 44321.7717800625

# Append new row 9 with the latest day's data
$ws.Cells.Item(9, 1).Value = 44322.77261097741
$ws.Cells.Item(9, 2).Value = 72481
$ws.Cells.Item(9, 3).Value = 60952
$ws.Cells.Item(9, 4).Value = 3222
$ws.Cells.Item(9, 5).Value = 1993
$ws.Cells.Item(9, 6).Value = 1409
$ws.Cells.Item(9, 7).Value = 19013
$ws.Cells.Item(9, 8).Value = 1359
$ws.Cells.Item(9, 9).Value = 811
$ws.Cells.Item(9, 10).Value = 210

# Make sure the new date cell uses the same number format as the rest of column A
$ws.Cells.Item(9, 1).NumberFormat = $ws.Cells.Item(8, 1).NumberFormat
